$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new column before column N (14th column), shifting existing
# "Late" / "Outstanding" columns one place to the right.
$ws.Columns("N:N").Insert()

# Give the freshly inserted column a plain default width (matches the
# author's Excel behaviour of a new, non bestFit column raw width of 10).
$ws.Columns("N:N").ColumnWidth = 9.140625

# Update the active selection to reflect where the author ended up editing.
$ws.Range("S8").Select()
